# update data for VopX and VopFull
# Sheet "CCmd" columns J (Vop-optimized) and L (EXPfull-optimized), rows 2-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CCmd")

$ws.Range("J2").Value = -1.1655
$ws.Range("L2").Value = -1.1634

$ws.Range("J3").Value = -1.1862
$ws.Range("L3").Value = -1.1845

$ws.Range("J4").Value = -1.2046
$ws.Range("L4").Value = -1.2032

$ws.Range("J5").Value = -1.22
$ws.Range("L5").Value = -1.219

$ws.Range("J6").Value = -1.232
$ws.Range("L6").Value = -1.2314

$ws.Range("J7").Value = -1.2399
$ws.Range("L7").Value = -1.2398

$ws.Range("J8").Value = -1.243
$ws.Range("L8").Value = -1.2434

$ws.Range("J9").Value = -1.2405
$ws.Range("L9").Value = -1.2414

$ws.Range("J10").Value = -1.2314
$ws.Range("L10").Value = -1.2329

$ws.Range("J11").Value = -1.2149
$ws.Range("L11").Value = -1.2169

$ws.Range("J12").Value = -1.1895
$ws.Range("L12").Value = -1.1921

$ws.Range("J13").Value = -1.1541
$ws.Range("L13").Value = -1.1572

$ws.Range("J14").Value = -1.107
$ws.Range("L14").Value = -1.1107

$ws.Range("J15").Value = -1.0466
$ws.Range("L15").Value = -1.0507

$ws.Range("J16").Value = -0.9708
$ws.Range("L16").Value = -0.9754

# Match the author's final selection state on the CCmd sheet.
$ws.Activate()
$ws.Range("L3").Select()
